$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.664.87'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '1.656.31'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.514'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.23'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.260'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0619'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0878'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.39%  '
$ws.Range('D12').Value = '1.895.20'
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').Value = '1.664.07'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.547'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '246.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.63%  '
$ws.Range('D18').Value = '27.716.32'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').Value = '0.0₃0727'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.39%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.54'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.63%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.111'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('E30').Value = '  +6.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0498'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('D33').Value = '1.427.84'
$ws.Range('E33').Value = '  -7.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.924'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.578'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.18%  '
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('E40').Value = '  -1.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '68.94'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.40'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.66%  '
$ws.Range('D45').Value = '1.803.64'
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.786'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.67'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').Value = '0.0₆0106'
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('E50').Value = '  -3.50%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0508'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.67%  '
